$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "2025/12/04 16:00"
$ws.Range("B47").Value = "-"
$ws.Range("C47").Value = "-"
$ws.Range("D47").Value = "-"
$ws.Range("E47").Value = "-"
$ws.Range("F47").Value = "-"
$ws.Range("G47").Value = "-"
